$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with a new "2021" column (O), mirroring the formatting
# already used for the neighbouring "2020" column (N) on each data row.
$ws.Range("N4:N14").Copy() | Out-Null
$ws.Range("O4:O14").PasteSpecial(-4122) | Out-Null
$ws.Range("N16:N17").Copy() | Out-Null
$ws.Range("O16:O17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 11.7
$ws.Range("O6").Value = 16.4
$ws.Range("O7").Value = 9.7
$ws.Range("O8").Value = 12.1
$ws.Range("O9").Value = 5.3
$ws.Range("O10").Value = 4.7
$ws.Range("O11").Value = 3.4
$ws.Range("O12").Value = 18.8
$ws.Range("O13").Value = 19.6
$ws.Range("O14").Value = 6.9
$ws.Range("O16").Value = 12.8
$ws.Range("O17").Value = 11

# Update the view: selection moves to R11 (this also drops the stale
# topLeftCell scroll position that used to keep column E in view).
$ws.Range("R11").Select() | Out-Null
